$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: add Q1's companion resistor R3 to the 10k resistor group designator list
$ws.Range("B10").Value = "R1, R2, R3, R5, R6, R7, R8, R9, R10, R11"

# New row 17: level shifter transistor (2N7002, Q1, SOT-23, LCSC C8545)
$ws.Range("B17").Value = "Q1"
$ws.Range("A17").Value = "2N7002"
$ws.Range("D17").Value = "C8545"
$ws.Range("C17").Value = "SOT-23"

# Row 15 (SRV05-4 / U2): strip "Package_TO_SOT_SMD:" library prefix from footprint
$ws.Range("C15").Value = "SOT-23-6"

# Row 13 (SW_Push / U-switches): strip "Button_Switch_SMD:" library prefix from footprint
$ws.Range("C13").Value = "SW_SPST_TL3342"

# Row 10: "10k" -> "10K"
$ws.Range("A10").Value = "10K"

# Row 16 (ESP32-S2-WROVER / U3): strip "RF_Module:" library prefix from footprint
$ws.Range("C16").Value = "ESP32-S2-WROVER"

$ws.Range("A10").Select()
